# "Generate Report for Handback"
#
# Adds "Latest Target File" (E) / "Latest Handback File" (F) data to the
# zh-cn and de-de report sheets, flips the Status text from "Ready for
# handoff" to "Handed back: in sync with en-US", and refreshes the
# "Latest Handback DateTime" (G) timestamps, for each localized-file row.

$wb = $excel.ActiveWorkbook

function Update-ReportSheet {
    param(
        $SheetName,
        $MdFile1,
        $XlfFile1,
        $XlfHandoffUrl1,
        $MdFile2,
        $XlfFile2,
        $XlfHandoffUrl2,
        $HandbackDateTime1,
        $HandbackDateTime2
    )

    $ws = $wb.Worksheets.Item($SheetName)

    $mdUrl1 = "https://github.com/OpenLocalizationTest/oltest/blob/1f633aa62bdbbcfdc3182fcfa5ebee83bf5800a7/e2e/$MdFile1"
    $mdUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/1f633aa62bdbbcfdc3182fcfa5ebee83bf5800a7/e2e/$MdFile2"
    $cfgUrl = "https://github.com/OpenLocalizationTest/oltest/blob/1f633aa62bdbbcfdc3182fcfa5ebee83bf5800a7/.localization-config"

    # Status text: handed back, now in sync with en-US.
    $ws.Range("B2").Value = "Handed back: in sync with en-US"
    $ws.Range("B3").Value = "Handed back: in sync with en-US"

    # New columns for this handback: Latest Target File / Latest Handback File.
    $ws.Range("E2").Value = $MdFile1
    $ws.Range("F2").Value = $XlfFile1
    $ws.Range("E3").Value = $MdFile2
    $ws.Range("F3").Value = $XlfFile2

    # Latest Handback DateTime moves forward to the handback run.
    $ws.Range("G2").Value = $HandbackDateTime1
    $ws.Range("G3").Value = $HandbackDateTime2

    # Rebuild every hyperlink on the sheet, in row order, so new E/F links
    # land right after their row's existing A/C links.
    $ws.Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl1, [Type]::Missing, [Type]::Missing, $MdFile1)
    $ws.Hyperlinks.Add($ws.Range("C2"), $XlfHandoffUrl1, [Type]::Missing, [Type]::Missing, $XlfFile1)
    $ws.Hyperlinks.Add($ws.Range("E2"), $mdUrl1, [Type]::Missing, [Type]::Missing, $MdFile1)
    $ws.Hyperlinks.Add($ws.Range("F2"), $XlfHandoffUrl1, [Type]::Missing, [Type]::Missing, $XlfFile1)

    $ws.Hyperlinks.Add($ws.Range("A3"), $mdUrl2, [Type]::Missing, [Type]::Missing, $MdFile2)
    $ws.Hyperlinks.Add($ws.Range("C3"), $XlfHandoffUrl2, [Type]::Missing, [Type]::Missing, $XlfFile2)
    $ws.Hyperlinks.Add($ws.Range("E3"), $mdUrl2, [Type]::Missing, [Type]::Missing, $MdFile2)
    $ws.Hyperlinks.Add($ws.Range("F3"), $XlfHandoffUrl2, [Type]::Missing, [Type]::Missing, $XlfFile2)

    $ws.Hyperlinks.Add($ws.Range("A4"), $cfgUrl, [Type]::Missing, [Type]::Missing, ".localization-config")

    $ws.Range("A2,C2,E2,F2,A3,C3,E3,F3,A4").Style = "HyperLink"
}

Update-ReportSheet `
    "zh-cn" `
    "91388e6b-1fec-4809-97e7-20e9b93df1be.md" `
    "91388e6b-1fec-4809-97e7-20e9b93df1be.480380839f5c1b35fb8c122dee73a5541a7799df.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d1d39d757907cf9b289d415e57752dd9971a3695/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/91388e6b-1fec-4809-97e7-20e9b93df1be.480380839f5c1b35fb8c122dee73a5541a7799df.zh-cn.xlf" `
    "af327491-d4ff-4eee-88ba-36f2b281dbfc.md" `
    "af327491-d4ff-4eee-88ba-36f2b281dbfc.049c8df917cc109999cf2faceb4dd68dc8e0754f.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d1d39d757907cf9b289d415e57752dd9971a3695/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/af327491-d4ff-4eee-88ba-36f2b281dbfc.049c8df917cc109999cf2faceb4dd68dc8e0754f.zh-cn.xlf" `
    "2016-02-06 04:32:21" `
    "2016-02-06 04:32:21"

Update-ReportSheet `
    "de-de" `
    "91388e6b-1fec-4809-97e7-20e9b93df1be.md" `
    "91388e6b-1fec-4809-97e7-20e9b93df1be.480380839f5c1b35fb8c122dee73a5541a7799df.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f4603ea9ab96a39a435be4a054504d861ea1e4ae/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/91388e6b-1fec-4809-97e7-20e9b93df1be.480380839f5c1b35fb8c122dee73a5541a7799df.de-de.xlf" `
    "af327491-d4ff-4eee-88ba-36f2b281dbfc.md" `
    "af327491-d4ff-4eee-88ba-36f2b281dbfc.049c8df917cc109999cf2faceb4dd68dc8e0754f.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f4603ea9ab96a39a435be4a054504d861ea1e4ae/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/af327491-d4ff-4eee-88ba-36f2b281dbfc.049c8df917cc109999cf2faceb4dd68dc8e0754f.de-de.xlf" `
    "2016-02-06 04:32:40" `
    "2016-02-06 04:32:40"

Write-Output "done"
